# Applies updated sensitivity values (Exact 1-Year Shifting No Double-Deployment - Total capex)
# to Sheet1 of the workbook: replaces the numeric data in D2:AH4 with the new values
# from the IEV model re-run (copy/paste shifting [2021:2050] -> [2020:2049]).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2" = 9996.1802109232158;
    "E2" = 10337.42880619616;
    "F2" = 11251.158026241977;
    "G2" = 9150.9433751602101;
    "H2" = 8604.7241095565896;
    "I2" = 6596.8199411359947;
    "J2" = 6638.5086766273898;
    "K2" = 6579.7987188229608;
    "L2" = 5446.2508487533078;
    "M2" = 5288.6939738602196;
    "N2" = 3902.5961006952762;
    "O2" = 3920.5532453150945;
    "P2" = 3955.9004194149011;
    "Q2" = 3997.6330012897502;
    "R2" = 4045.6948583918588;
    "S2" = 4108.6213687890431;
    "T2" = 4179.3847687297166;
    "U2" = 4258.9683020458515;
    "V2" = 4601.2992917750844;
    "W2" = 5102.8241745830101;
    "X2" = 5732.3727324716956;
    "Y2" = 6046.7377561118356;
    "Z2" = 6370.3112347253063;
    "AA2" = 6703.6346798791756;
    "AB2" = 7047.2906548054107;
    "AC2" = 7401.9010865484015;
    "AD2" = 7768.1226601829667;
    "AE2" = 8146.6380118330471;
    "AF2" = 7771.4310750614522;
    "AG2" = 7595.9765851570082;
    "AH2" = 4764.4463864408744;
    "D3" = 9996.1802109232121;
    "E3" = 10363.239393158467;
    "F3" = 11297.455097262868;
    "G3" = 9221.546361714969;
    "H3" = 8702.7847084809673;
    "I3" = 6687.8339961961537;
    "J3" = 6750.4805629987914;
    "K3" = 6717.442752471723;
    "L3" = 6119.1662826357069;
    "M3" = 6729.8893007187589;
    "N3" = 6236.8405399262319;
    "O3" = 6432.4507625996466;
    "P3" = 6653.9049890881997;
    "Q3" = 6890.1837934022878;
    "R3" = 7141.2650343687537;
    "S3" = 7413.4115157674669;
    "T3" = 7701.5591662693605;
    "U3" = 8006.887671842881;
    "V3" = 8583.5916314088736;
    "W3" = 9328.4403337759541;
    "X3" = 10210.801041384675;
    "Y3" = 10787.653170686113;
    "Z3" = 11383.728116110904;
    "AA3" = 11999.851215137627;
    "AB3" = 12636.799221540388;
    "AC3" = 13295.261386109725;
    "AD3" = 13975.792837384814;
    "AE3" = 14678.756941928528;
    "AF3" = 12818.305296404124;
    "AG3" = 11698.525590682984;
    "AH3" = 4453.9159480525213;
    "D4" = 9996.1802109232121;
    "E4" = 10424.736870822924;
    "F4" = 11407.764603017042;
    "G4" = 9389.7682465964008;
    "H4" = 8936.4283486245458;
    "I4" = 6904.6882206424525;
    "J4" = 7017.2698689231056;
    "K4" = 7045.399683374666;
    "L4" = 6924.3301993856103;
    "M4" = 8222.7149135102136;
    "N4" = 8506.1444026307418;
    "O4" = 8973.189847151385;
    "P4" = 9484.3890921076054;
    "Q4" = 10029.277795761949;
    "R4" = 10608.583766017713;
    "S4" = 11223.914952963454;
    "T4" = 11875.749858725238;
    "U4" = 12566.717851763658;
    "V4" = 13552.764960748795;
    "W4" = 14732.693851302887;
    "X4" = 16078.372195807153;
    "Y4" = 17149.063232286408;
    "Z4" = 18272.136209689084;
    "AA4" = 19451.155591229035;
    "AB4" = 20689.683476901268;
    "AC4" = 21991.180237662295;
    "AD4" = 23358.888152875101;
    "AE4" = 24795.689675839356;
    "AF4" = 20831.493825390953;
    "AG4" = 18358.464145661434;
    "AH4" = 3977.144719262611;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
